# Apply cryptos price/volume update (Mon Jul  8 09:34:28 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain numeric-looking string (e.g. "140.02") need a
# leading quote-prefix so Excel stores them as TEXT (matching the original inline-string
# column) instead of auto-converting them to a number.
$ws.Range("D2").Value = "57.659.11"
$ws.Range("E2").Value = "  -0.04%  "
$ws.Range("D3").Value = "3.061.93"
$ws.Range("E3").Value = "  +1.56%  "
$ws.Range("E4").Value = "  +0.61%  "
$ws.Range("D5").Value = "'516.13"
$ws.Range("E5").Value = "  +0.79%  "
$ws.Range("D6").Value = "'140.02"
$ws.Range("E6").Value = "  -0.06%  "
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.434"
$ws.Range("E8").Value = "  -1.01%  "
$ws.Range("D9").Value = "'7.28"
$ws.Range("E9").Value = "  -3.83%  "
$ws.Range("D10").Value = "'0.109"
$ws.Range("E10").Value = "  -1.35%  "
$ws.Range("D11").Value = "'0.371"
$ws.Range("E11").Value = "  +1.27%  "
$ws.Range("D12").Value = "3.576.81"
$ws.Range("E12").Value = "  +1.40%  "
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "'26.80"
$ws.Range("E14").Value = "  +0.80%  "
$ws.Range("D15").Value = "'0.0000167"
$ws.Range("E15").Value = "  +1.33%  "
$ws.Range("D16").Value = "57.638.82"
$ws.Range("E16").Value = "  +0.04%  "
$ws.Range("D17").Value = "'6.22"
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "3.065.97"
$ws.Range("E18").Value = "  +1.75%  "
$ws.Range("D19").Value = "'13.34"
$ws.Range("E19").Value = "  +3.57%  "
$ws.Range("D20").Value = "'8.18"
$ws.Range("E20").Value = "  +2.21%  "
$ws.Range("D21").Value = "'329.41"
$ws.Range("E21").Value = "  -0.81%  "
$ws.Range("E22").Value = "  +0.21%  "
$ws.Range("D23").Value = "'0.508"
$ws.Range("E23").Value = "  +1.52%  "
$ws.Range("D24").Value = "'65.38"
$ws.Range("E24").Value = "  +1.16%  "
$ws.Range("E25").Value = "  +1.34%  "
$ws.Range("D26").Value = "'0.165"
$ws.Range("E26").Value = "  -2.58%  "
$ws.Range("D27").Value = "'0.994"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").Value = "0.0₃0900"
$ws.Range("E28").Value = "  -3.04%  "
$ws.Range("D29").Value = "'6.68"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "'7.27"
$ws.Range("E30").Value = "  -2.00%  "
$ws.Range("D31").Value = "'1.81"
$ws.Range("E31").Value = "  -0.60%  "
$ws.Range("D32").Value = "'1.21"
$ws.Range("E32").Value = "  +1.66%  "
$ws.Range("D33").Value = "'20.80"
$ws.Range("E33").Value = "  +0.63%  "
$ws.Range("D34").Value = "'154.06"
$ws.Range("E34").Value = "  -0.41%  "
$ws.Range("D35").Value = "'4.62"
$ws.Range("E35").Value = "  -2.70%  "
$ws.Range("D36").Value = "'5.86"
$ws.Range("E36").Value = "  -0.45%  "
$ws.Range("D37").Value = "'25.50"
$ws.Range("E37").Value = "  +3.49%  "
$ws.Range("D38").Value = "'1.27"
$ws.Range("E38").Value = "  -1.09%  "
$ws.Range("D39").Value = "'0.0677"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'37.12"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("D41").Value = "'3.88"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("E42").Value = "  +0.43%  "
$ws.Range("E43").Value = "  +2.96%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "2.205.53"
$ws.Range("E44").Value = "  -1.66%  "
$ws.Range("B45").Value = "Stacks"
$ws.Range("C45").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D45").Value = "'1.39"
$ws.Range("E45").Value = "  -1.57%  "
$ws.Range("D46").Value = "'6.11"
$ws.Range("E46").Value = "  +1.06%  "
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").Value = "'0.0246"
$ws.Range("E47").Value = "  +2.11%  "
$ws.Range("B48").Value = "ONDO"
$ws.Range("C48").Value = "https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo"
$ws.Range("D48").Value = "'0.952"
$ws.Range("E48").Value = "  -3.70%  "
$ws.Range("D49").Value = "'19.95"
$ws.Range("E49").Value = "  +2.11%  "
$ws.Range("B50").Value = "TheGraph"
$ws.Range("C50").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D50").Value = "'0.184"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("B51").Value = "Notcoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/2L2Y4ghjj+notcoin-not"
$ws.Range("D51").Value = "'0.0171"
$ws.Range("E51").Value = "  +8.91%  "
